$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.851.56"
$ws.Range("E2").Value = "'  -2.59%  "
$ws.Range("D3").Value = "'1.966.81"
$ws.Range("E3").Value = "'  -1.72%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "'  -0.21%  "
$ws.Range("D5").Value = "'324.12"
$ws.Range("E5").Value = "'  -2.13%  "
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "'  -0.07%  "
$ws.Range("D7").Value = "'0.4773"
$ws.Range("E7").Value = "'  -4.55%  "
$ws.Range("D8").Value = "'0.4045"
$ws.Range("E8").Value = "'  -4.68%  "
$ws.Range("D9").Value = "'53.99"
$ws.Range("E9").Value = "'  -0.93%  "
$ws.Range("D10").Value = "'0.08548"
$ws.Range("E10").Value = "'  -5.76%  "
$ws.Range("D11").Value = "'1.062"
$ws.Range("E11").Value = "'  -5.32%  "
$ws.Range("D12").Value = "'22.49"
$ws.Range("E12").Value = "'  -4.12%  "
$ws.Range("D13").Value = "'1.964.43"
$ws.Range("E13").Value = "'  -5.51%  "
$ws.Range("D14").Value = "'7.671"
$ws.Range("E14").Value = "'  -5.26%  "
$ws.Range("D15").Value = "'6.257"
$ws.Range("E15").Value = "'  -3.80%  "
$ws.Range("D16").Value = "'1.014"
$ws.Range("E16").Value = "'  +0.03%  "
$ws.Range("D17").Value = "'90.17"
$ws.Range("E17").Value = "'  -4.77%  "
$ws.Range("D18").Value = "'0.00001071"
$ws.Range("E18").Value = "'  -4.01%  "
$ws.Range("D19").Value = "'0.06624"
$ws.Range("E19").Value = "'  -0.61%  "
$ws.Range("D20").Value = "'18.71"
$ws.Range("E20").Value = "'  -5.68%  "
$ws.Range("D21").Value = "'1.012"
$ws.Range("E21").Value = "'  -0.18%  "
$ws.Range("D22").Value = "'5.796"
$ws.Range("E22").Value = "'  -2.97%  "
$ws.Range("D23").Value = "'28.862.30"
$ws.Range("E23").Value = "'  -2.56%  "
$ws.Range("D24").Value = "'11.58"
$ws.Range("E24").Value = "'  -3.79%  "
$ws.Range("D25").Value = "'2.295"
$ws.Range("E25").Value = "'  +0.51%  "
$ws.Range("D26").Value = "'2.242.49"
$ws.Range("E26").Value = "'  -2.65%  "
$ws.Range("D27").Value = "'154.21"
$ws.Range("E27").Value = "'  -2.85%  "
$ws.Range("D28").Value = "'20.26"
$ws.Range("E28").Value = "'  -2.26%  "
$ws.Range("D29").Value = "'5.974"
$ws.Range("E29").Value = "'  -6.57%  "
$ws.Range("D30").Value = "'2.150"
$ws.Range("E30").Value = "'  -6.97%  "
$ws.Range("D31").Value = "'124.41"
$ws.Range("E31").Value = "'  -3.27%  "
$ws.Range("D32").Value = "'1.008"
$ws.Range("E32").Value = "'  -4.71%  "
$ws.Range("D33").Value = "'0.09643"
$ws.Range("E33").Value = "'  -3.00%  "
$ws.Range("D34").Value = "'1.464"
$ws.Range("E34").Value = "'  -6.94%  "
$ws.Range("D35").Value = "'5.691"
$ws.Range("E35").Value = "'  -2.76%  "
$ws.Range("D36").Value = "'3.693"
$ws.Range("E36").Value = "'  -2.52%  "
$ws.Range("D37").Value = "'0.02347"
$ws.Range("E37").Value = "'  -5.04%  "
$ws.Range("D38").Value = "'0.06217"
$ws.Range("E38").Value = "'  -2.40%  "
$ws.Range("D39").Value = "'1.267"
$ws.Range("E39").Value = "'  -3.45%  "
$ws.Range("D40").Value = "'8.775"
$ws.Range("E40").Value = "'  -7.73%  "
$ws.Range("D41").Value = "'0.6253"
$ws.Range("E41").Value = "'  -5.00%  "
$ws.Range("D42").Value = "'11.09"
$ws.Range("E42").Value = "'  -5.35%  "
$ws.Range("D43").Value = "'1.011"
$ws.Range("E43").Value = "'  -0.07%  "
$ws.Range("D44").Value = "'0.1924"
$ws.Range("E44").Value = "'  -6.68%  "
$ws.Range("D45").Value = "'1.354"
$ws.Range("E45").Value = "'  +5.17%  "
$ws.Range("D46").Value = "'0.5985"
$ws.Range("E46").Value = "'  -5.80%  "
$ws.Range("D47").Value = "'12.94"
$ws.Range("E47").Value = "'  -4.83%  "
$ws.Range("D48").Value = "'2.088"
$ws.Range("E48").Value = "'  -5.54%  "
$ws.Range("D49").Value = "'3.429"
$ws.Range("E49").Value = "'  -2.88%  "
$ws.Range("D50").Value = "'0.00000000336"
$ws.Range("E50").Value = "'  +0.06%  "
$ws.Range("D51").Value = "'2.113"
$ws.Range("E51").Value = "'  +5.84%  "

Write-Output "Updated cryptos list"